$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the row-label strings in column A (rows 2-6).
$ws.Range("A2").Value = "attendence"
$ws.Range("A3").Value = "materials"
$ws.Range("A4").Value = "sessions"
$ws.Range("A5").Value = "total_Engagement"
$ws.Range("A6").Value = "total_Cost"

# Add a new, empty cell at A9 that carries the wrap-text format A4 used
# to have, by copying A4's (still wrap-text) formatting onto it.
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A4 itself loses the wrap-text formatting.
$ws.Range("A4").WrapText = $false
